# Daily attendance processing - normalize the "Recorded By" (column G) entries.
# Each cell holds a comma-separated list of recorders; re-order the list
# according to a canonical recorder priority instead of the arbitrary order
# they were appended in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RecorderPriority($name) {
    if ($name.CompareTo("backup@backdoor.com") -eq 0) { return 0 }
    if ($name.CompareTo("dnasr281@gmail.com") -eq 0) { return 1 }
    if ($name.CompareTo("admin@admin.com") -eq 0) { return 2 }
    if ($name.CompareTo("System") -eq 0) { return 3 }
    if ($name.CompareTo("system") -eq 0) { return 4 }
    return 50
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Text

    if ($txt -eq $null -or $txt -eq "") {
        continue
    }

    $parts = $txt -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    $items = @()
    foreach ($p in $parts) {
        $pri = Get-RecorderPriority $p
        $items += [PSCustomObject]@{ Name = $p; Pri = $pri }
    }

    $sortedItems = $items | Sort-Object -Property Pri

    $names = @()
    foreach ($o in $sortedItems) {
        $names += $o.Name
    }

    $newTxt = $names -join ", "

    if ($newTxt -cne $txt) {
        $cell.Value = $newTxt
    }
}
